$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.802.58"
$ws.Cells.Item(2, 5).Value = "  +0.68%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.467.23"
$ws.Cells.Item(3, 5).Value = "  +0.53%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "560.79"
$ws.Cells.Item(5, 5).Value = "  +0.34%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "163.16"
$ws.Cells.Item(6, 5).Value = "  -0.10%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.06%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.513"
$ws.Cells.Item(8, 5).Value = "  +1.97%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.155"
$ws.Cells.Item(9, 5).Value = "  +2.89%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.70%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.333"
$ws.Cells.Item(11, 5).Value = "  -1.47%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.86"
$ws.Cells.Item(12, 5).Value = "  +1.25%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "68.645.54"
$ws.Cells.Item(13, 5).Value = "  +0.57%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000169"
$ws.Cells.Item(14, 5).Value = "  -0.49%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "23.60"
$ws.Cells.Item(15, 5).Value = "  +1.28%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "10.68"
$ws.Cells.Item(16, 5).Value = "  -2.72%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "336.71"
$ws.Cells.Item(17, 5).Value = "  -1.81%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.93"
$ws.Cells.Item(18, 5).Value = "  -3.60%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.79"
$ws.Cells.Item(19, 5).Value = "  +0.26%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "SuiNetwork"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.89"
$ws.Cells.Item(20, 5).Value = "  +1.85%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "Dai"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.00"
$ws.Cells.Item(21, 5).Value = "  +0.01%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "66.48"
$ws.Cells.Item(22, 5).Value = "  -1.65%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.67"
$ws.Cells.Item(23, 5).Value = "  -1.13%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.29"
$ws.Cells.Item(24, 5).Value = "  +1.67%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "0.0₃0825"
$ws.Cells.Item(25, 5).Value = "  -1.03%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.23"
$ws.Cells.Item(26, 5).Value = "  -0.37%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.02%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "431.16"
$ws.Cells.Item(28, 5).Value = "  -0.57%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.15"
$ws.Cells.Item(29, 5).Value = "  -1.98%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.62"
$ws.Cells.Item(30, 5).Value = "  -2.82%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "159.27"
$ws.Cells.Item(31, 5).Value = "  +1.28%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "19.02"
$ws.Cells.Item(32, 5).Value = "  +0.11%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.02%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -1.19%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "17.83"
$ws.Cells.Item(35, 5).Value = "  -0.10%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "RenderToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.44"
$ws.Cells.Item(36, 5).Value = "  -0.13%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.300"
$ws.Cells.Item(37, 5).Value = "  -1.79%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.47"
$ws.Cells.Item(38, 5).Value = "  -2.89%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.08"
$ws.Cells.Item(39, 5).Value = "  -0.69%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.08"
$ws.Cells.Item(40, 5).Value = "  +0.94%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.36"
$ws.Cells.Item(41, 5).Value = "  +0.38%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "129.98"
$ws.Cells.Item(42, 5).Value = "  -2.79%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0719"
$ws.Cells.Item(43, 5).Value = "  +0.24%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.484"
$ws.Cells.Item(44, 5).Value = "  +0.04%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.564"
$ws.Cells.Item(45, 5).Value = "  +0.56%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0918"
$ws.Cells.Item(46, 5).Value = "  +1.12%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.33%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.39"
$ws.Cells.Item(48, 5).Value = "  -2.48%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "5.01"
$ws.Cells.Item(49, 5).Value = "  -6.60%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "16.91"
$ws.Cells.Item(50, 5).Value = "  -3.72%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "0.0₆0207"
$ws.Cells.Item(51, 5).Value = "  +0.67%  "
